$wb = $excel.ActiveWorkbook

# --- eval_metrics sheet: relabel rows 3-4 (swap) and 7-13 (permute) ---
$ws = $wb.Worksheets.Item("eval_metrics")
$ws.Range("A3").Value = "MAE"
$ws.Range("A4").Value = "MAPE"
$ws.Range("A7").Value = "nRMSE"
$ws.Range("A8").Value = "RMSPE"
$ws.Range("A9").Value = "accuracy"
$ws.Range("A10").Value = "IP"
$ws.Range("A11").Value = "RMSLE"
$ws.Range("A12").Value = "precision"
$ws.Range("A13").Value = "RMdSPE"

# --- limitation_cats sheet: swap rows 3-4 ---
$ws = $wb.Worksheets.Item("limitation_cats")
$ws.Range("A3").Value = "data"
$ws.Range("A4").Value = "note factors that were not accounted for in their model"

# --- data_cats sheet: swap rows 6-7 and rows 8-9 ---
$ws = $wb.Worksheets.Item("data_cats")
$ws.Range("A6").Value = "climate"
$ws.Range("A7").Value = "policy"
$ws.Range("A8").Value = "demographics"
$ws.Range("A9").Value = "mobility"

# --- journal_subjects sheet: rotate rows 2-4 ---
$ws = $wb.Worksheets.Item("journal_subjects")
$ws.Range("A2").Value = "Statistical and Nonlinear Physics"
$ws.Range("A3").Value = "General Physics and Astronomy"
$ws.Range("A4").Value = "General Mathematics"

# --- region_level sheet: swap rows 3-4 ---
$ws = $wb.Worksheets.Item("region_level")
$ws.Range("A3").Value = "county or smaller"
$ws.Range("A4").Value = "state"

# --- target_cats sheet: permute rows 6, 8-13 ---
$ws = $wb.Worksheets.Item("target_cats")
$ws.Range("A6").Value = "end dates of pandemic"
$ws.Range("A8").Value = "peak cases date"
$ws.Range("A9").Value = "ventilators"
$ws.Range("A10").Value = "symptomatic cases"
$ws.Range("A11").Value = "peak deaths"
$ws.Range("A12").Value = "total deaths"
$ws.Range("A13").Value = "attack rate"
